$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row labels: abbreviate "Current" -> "Crnt" in the test names (A2:A5) ---
$ws.Range("A2").Value = "Test 1 - High MFR Low Crnt"
$ws.Range("A3").Value = "Test 2 - High MFR High Crnt"
$ws.Range("A4").Value = "Test 3 - Low MFR High Crnt"
$ws.Range("A5").Value = "Test 4 - Low MFR Low Crnt"

# --- Inlet / Outlet temperature data: corrected (Celsius -> Kelvin offset fix, +0.15) ---
$ws.Range("G2").Value = 310.65
$ws.Range("H2").Value = 325.15
$ws.Range("G3").Value = 312.65
$ws.Range("H3").Value = 340.15
$ws.Range("G4").Value = 313.95
$ws.Range("H4").Value = 347.15
$ws.Range("G5").Value = 314.32
$ws.Range("H5").Value = 332.15

# --- Selection moved to H3 (where the author was last working) ---
$ws.Range("H3").Select() | Out-Null
